# Fixed a bug in Bomb
# The data rows (2-25, columns A-F) on the active sheet were shuffled into a
# new order. Capture the existing values first, then write them back out in
# their new row positions so no data is lost while permuting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of source row -> destination row for the permutation.
$rowMap = @{
    2  = 7
    3  = 9
    4  = 2
    5  = 10
    6  = 12
    7  = 15
    8  = 14
    9  = 6
    10 = 11
    11 = 3
    12 = 5
    13 = 8
    14 = 13
    15 = 4
    16 = 20
    17 = 16
    18 = 17
    19 = 18
    20 = 19
    21 = 21
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

$firstCol = 1  # A
$lastCol  = 6  # F

# Snapshot all the current values for rows 2-25, columns A-F before writing
# anything back, since several rows both read from and write to overlapping
# row numbers.
$snapshot = @{}
foreach ($srcRow in $rowMap.Keys) {
    $rowValues = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value2
    }
    $snapshot[$srcRow] = $rowValues
}

foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    $rowValues = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($dstRow, $col).Value2 = $rowValues[$col]
    }
}
